$d = $word.ActiveDocument

$replacements = @(
    @{old="19×16="; new="87×12="},
    @{old="54×77="; new="98×48="},
    @{old="40×77="; new="47×79="},
    @{old="49×25="; new="14×42="},
    @{old="30×29="; new="67×93="},
    @{old="44×11="; new="19×21="},
    @{old="42×83="; new="90×73="},
    @{old="18×52="; new="21×88="},
    @{old="56×49="; new="43×76="},
    @{old="16×19="; new="42×41="},
    @{old="87×43="; new="68×82="},
    @{old="32×39="; new="21×68="},
    @{old="76×20="; new="83×21="},
    @{old="60×69="; new="23×49="},
    @{old="41×45="; new="99×69="},
    @{old="98×68="; new="98×34="},
    @{old="22×88="; new="60×22="},
    @{old="55×16="; new="59×69="},
    @{old="46×69="; new="97×34="},
    @{old="49×95="; new="71×13="},
    @{old="63×19="; new="14×61="},
    @{old="59×20="; new="39×46="},
    @{old="84×86="; new="33×44="},
    @{old="89×99="; new="42×92="},
    @{old="87×22="; new="71×22="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
